# Auto update: 2025-12-01 01:09:22
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update date column (A2:A5): 2025-11-29 -> 2025-12-01
# Force text entry (avoid Excel auto-converting the date-like string into
# a date serial number), then drop the temporary "@" number format so the
# cells end up back on the sheet's default style, matching the source cells.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2:A5").Value = "2025-12-01"
$ws.Range("A2:A5").NumberFormat = "General"
$ws.Range("A2:A5").ClearFormats()

# Update MACRO_SCORE column (N2:N5) with refreshed value
$ws.Range("N2:N5").Value = 85.87246918135976
